$d = $word.ActiveDocument

# -------------------------------------------------------------------
# 1) Merge the 3 split runs that make up the stackoverflow URL into a
#    single run (keeping the Hyperlink character style).
# -------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("https://stackoverflow.com/questions/38512485/highlight-specific-points-in-matplotlib-scatterplot", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Text = "TEMP_MERGE_PLACEHOLDER_38512485"

    $rng2 = $d.Content
    $found2 = $rng2.Find.Execute("TEMP_MERGE_PLACEHOLDER_38512485", `
        $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $rng2.Text = "https://stackoverflow.com/questions/38512485/highlight-specific-points-in-matplotlib-scatterplot"
        $rng2.Style = "Hyperlink"
    }
}

# -------------------------------------------------------------------
# 2) Append the new content at the end of the document, before the
#    final (already existing) empty paragraph stays untouched -- we
#    add our new paragraphs right after it.
# -------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertAfter([char]13 + "Combinar columnas concretas de excel con pandas:" + [char]13 + "PLACEHOLDER_LINK_1" + [char]13 + "PLACEHOLDER_LINK_2" + [char]13 + [char]13)

# -------------------------------------------------------------------
# 3) Turn the two placeholders into real hyperlinks.
# -------------------------------------------------------------------
$linkRange1 = $d.Content
$f1 = $linkRange1.Find.Execute("PLACEHOLDER_LINK_1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f1) {
    $d.Hyperlinks.Add($linkRange1, "https://www.geeksforgeeks.org/how-to-merge-two-csv-files-by-specific-column-using-pandas-in-python/", "", "", "https://www.geeksforgeeks.org/how-to-merge-two-csv-files-by-specific-column-using-pandas-in-python/") | Out-Null
}

$restyleRange1 = $d.Content
$fr1 = $restyleRange1.Find.Execute("https://www.geeksforgeeks.org/how-to-merge-two-csv-files-by-specific-column-using-pandas-in-python/", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($fr1) {
    $restyleRange1.Style = "Hyperlink"
}

$linkRange2 = $d.Content
$f2 = $linkRange2.Find.Execute("PLACEHOLDER_LINK_2", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($f2) {
    $d.Hyperlinks.Add($linkRange2, "https://www.youtube.com/watch?v=dcQs8k9WGbY", "", "", "https://www.youtube.com/watch?v=dcQs8k9WGbY") | Out-Null
}

$restyleRange2 = $d.Content
$fr2 = $restyleRange2.Find.Execute("https://www.youtube.com/watch?v=dcQs8k9WGbY", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($fr2) {
    $restyleRange2.Style = "Hyperlink"
}

Write-Host "Edit complete"
